# Scheduled market-data refresh for the Leve profit tables.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (columns H/I/J) and the resulting LevePriceNQ/HQ + LeveProfitNQ/HQ
# (columns K/L/M/N) pulled from the market-board API for each Leve sheet.
$wb = $excel.ActiveWorkbook

# ALC row 11: Gotta Bounce
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4261.5
$ws.Range("I11").Value = 4261.5
$ws.Range("K11").Value = 4261.5
$ws.Range("M11").Value = -4121.5

# ALC row 129: Practical Command
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 730.3125
$ws.Range("I129").Value = 397.33334
$ws.Range("K129").Value = 1192.00002
$ws.Range("M129").Value = 3807.99998

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3791.6333
$ws.Range("I132").Value = 4101.815
$ws.Range("K132").Value = 12305.445
$ws.Range("M132").Value = -9775.445

# ARM row 2: Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1550834
$ws.Range("I2").Value = 2392.6
$ws.Range("J2").Value = 2103848.8
$ws.Range("K2").Value = 2392.6
$ws.Range("L2").Value = 2103848.8
$ws.Range("M2").Value = -2279.6
$ws.Range("N2").Value = -2104074.8

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6328.25
$ws.Range("I32").Value = 5348.705
$ws.Range("J32").Value = 24939.6
$ws.Range("K32").Value = 5348.705
$ws.Range("L32").Value = 24939.6
$ws.Range("M32").Value = -5061.705
$ws.Range("N32").Value = -25513.6

# ARM row 37: Get Shirty
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8479
$ws.Range("J37").Value = 9891.4
$ws.Range("L37").Value = 9891.4
$ws.Range("N37").Value = -10437.4

# ARM row 56: Feasting the Night Away
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 55000
$ws.Range("J56").Value = 55000
$ws.Range("L56").Value = 55000
$ws.Range("N56").Value = -56484

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2932.423
$ws.Range("I74").Value = 1126.3334
$ws.Range("J74").Value = 4480.5
$ws.Range("K74").Value = 1126.3334
$ws.Range("L74").Value = 4480.5
$ws.Range("M74").Value = -252.3334
$ws.Range("N74").Value = -6228.5

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2932.423
$ws.Range("I77").Value = 1126.3334
$ws.Range("J77").Value = 4480.5
$ws.Range("K77").Value = 5631.666999999999
$ws.Range("L77").Value = 22402.5
$ws.Range("M77").Value = -1263.666999999999
$ws.Range("N77").Value = -31138.5

# ARM row 116: No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1550834
$ws.Range("I116").Value = 2392.6
$ws.Range("J116").Value = 2103848.8
$ws.Range("K116").Value = 2392.6
$ws.Range("L116").Value = 2103848.8
$ws.Range("M116").Value = -98.59999999999991
$ws.Range("N116").Value = -2108436.8

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 24666.822
$ws.Range("I132").Value = 2411.625
$ws.Range("J132").Value = 202708.4
$ws.Range("K132").Value = 7234.875
$ws.Range("L132").Value = 608125.2
$ws.Range("M132").Value = -4704.875
$ws.Range("N132").Value = -613185.2

# BSM row 3: Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1550834
$ws.Range("I3").Value = 2392.6
$ws.Range("J3").Value = 2103848.8
$ws.Range("K3").Value = 2392.6
$ws.Range("L3").Value = 2103848.8
$ws.Range("M3").Value = -2278.6
$ws.Range("N3").Value = -2104076.8

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 80308.92999999999
$ws.Range("I134").Value = 111673.5
$ws.Range("J134").Value = 1897.5
$ws.Range("K134").Value = 335020.5
$ws.Range("L134").Value = 5692.5
$ws.Range("M134").Value = -332485.5
$ws.Range("N134").Value = -10762.5

# CRP row 50: The Arsenal of Theocracy
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10633.2
$ws.Range("J50").Value = 10633.2
$ws.Range("L50").Value = 10633.2
$ws.Range("N50").Value = -11883.2

# CRP row 59: Bow Down to Magic
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 13408.1
$ws.Range("I59").Value = 7666.6665
$ws.Range("J59").Value = 15868.714
$ws.Range("K59").Value = 7666.6665
$ws.Range("L59").Value = 15868.714
$ws.Range("M59").Value = -6521.6665
$ws.Range("N59").Value = -18158.714

# CRP row 60: Bowing to Greater Power
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 7207
$ws.Range("I60").Value = 2600
$ws.Range("J60").Value = 9301.091
$ws.Range("K60").Value = 2600
$ws.Range("L60").Value = 9301.091
$ws.Range("N60").Value = -10323.091
$ws.Range("M60").Value = -2089

# CRP row 62: Splinter in the Sewers
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2457.2
$ws.Range("I62").Value = 2125.2942
$ws.Range("J62").Value = 3162.5
$ws.Range("K62").Value = 2125.2942
$ws.Range("L62").Value = 3162.5
$ws.Range("M62").Value = -1501.2942
$ws.Range("N62").Value = -4410.5

# CRP row 65: The Lumber of Their Discontent (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2457.2
$ws.Range("I65").Value = 2125.2942
$ws.Range("J65").Value = 3162.5
$ws.Range("K65").Value = 10626.471
$ws.Range("L65").Value = 15812.5
$ws.Range("M65").Value = -7506.471
$ws.Range("N65").Value = -22052.5

# CRP row 68: Do You Even String Bow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 18805.5
$ws.Range("J68").Value = 18805.5
$ws.Range("L68").Value = 18805.5
$ws.Range("N68").Value = -20303.5

# CRP row 71: Win One Bow, Get Three Free (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 18805.5
$ws.Range("J71").Value = 18805.5
$ws.Range("L71").Value = 56416.5
$ws.Range("N71").Value = -63904.5

# CRP row 74: License to Heal
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 15242.75
$ws.Range("J74").Value = 16848.857
$ws.Range("L74").Value = 16848.857
$ws.Range("N74").Value = -18596.857

# CRP row 77: Purified Polyrhythm (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 15242.75
$ws.Range("J77").Value = 16848.857
$ws.Range("L77").Value = 50546.571
$ws.Range("N77").Value = -59282.571

# CRP row 99: O Pine
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 67244.39
$ws.Range("I99").Value = 33309.312
$ws.Range("J99").Value = 144810.28
$ws.Range("K99").Value = 33309.312
$ws.Range("L99").Value = 144810.28
$ws.Range("M99").Value = -31811.312
$ws.Range("N99").Value = -147806.28

# CRP row 126: A Better Conductor
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 67244.39
$ws.Range("I126").Value = 33309.312
$ws.Range("J126").Value = 144810.28
$ws.Range("K126").Value = 99927.93599999999
$ws.Range("L126").Value = 434430.84
$ws.Range("M126").Value = -97457.93599999999
$ws.Range("N126").Value = -439370.84

# CUL row 9: Jack of All Plates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1332.2222
$ws.Range("I9").Value = 750
$ws.Range("J9").Value = 1498.5714
$ws.Range("K9").Value = 2250
$ws.Range("L9").Value = 4495.7142
$ws.Range("M9").Value = -2026
$ws.Range("N9").Value = -4943.7142

# CUL row 14: Keep Your Powder Dry
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 113.73684
$ws.Range("I14").Value = 113.73684
$ws.Range("K14").Value = 341.21052
$ws.Range("M14").Value = -168.21052

# CUL row 131: The Mountain Steeped
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2553.638
$ws.Range("J131").Value = 1759.6545
$ws.Range("L131").Value = 5278.9635
$ws.Range("N131").Value = -15358.9635

# CUL row 132: More Mezcal
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 38462492
$ws.Range("I132").Value = 52632400
$ws.Range("J132").Value = 1306.1428
$ws.Range("K132").Value = 473691600
$ws.Range("L132").Value = 11755.2852
$ws.Range("M132").Value = -473689070
$ws.Range("N132").Value = -16815.2852

# GSM row 12: Horn of Plenty
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3069.9312
$ws.Range("I12").Value = 1336.4445
$ws.Range("J12").Value = 3850
$ws.Range("K12").Value = 1336.4445
$ws.Range("L12").Value = 3850
$ws.Range("M12").Value = -1196.4445
$ws.Range("N12").Value = -4130

# GSM row 70: Sky Is the Limit
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7582.5938
$ws.Range("I70").Value = 10015.471
$ws.Range("K70").Value = 10015.471
$ws.Range("M70").Value = -9745.471

# GSM row 73: Hulls of Broken Dreams (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7582.5938
$ws.Range("I73").Value = 10015.471
$ws.Range("K73").Value = 10015.471
$ws.Range("M73").Value = -9079.471

# GSM row 80: Needs More Prayerbell
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2893.5
$ws.Range("I80").Value = 3191
$ws.Range("J80").Value = 2596
$ws.Range("K80").Value = 3191
$ws.Range("L80").Value = 2596
$ws.Range("M80").Value = -2193
$ws.Range("N80").Value = -4592

# GSM row 83: With a Noise That Reaches Heaven (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2893.5
$ws.Range("I83").Value = 3191
$ws.Range("J83").Value = 2596
$ws.Range("K83").Value = 15955
$ws.Range("L83").Value = 12980
$ws.Range("M83").Value = -10963
$ws.Range("N83").Value = -22964

# LTW row 18: Simply the Best
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 8000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1670.125
$ws.Range("I136").Value = 949.15
$ws.Range("J136").Value = 5275
$ws.Range("K136").Value = 2847.45
$ws.Range("L136").Value = 15825
$ws.Range("M136").Value = -297.4499999999998
$ws.Range("N136").Value = -20925

# WVR row 122: Heavy Armoire
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1850.1628
$ws.Range("I122").Value = 1456.3
$ws.Range("J122").Value = 2759.077
$ws.Range("K122").Value = 4368.9
$ws.Range("L122").Value = 8277.231
$ws.Range("M122").Value = -1918.9
$ws.Range("N122").Value = -13177.231

# WVR row 124: Hot Heads
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 49999.332
$ws.Range("J124").Value = 49999.332
$ws.Range("L124").Value = 49999.332
$ws.Range("N124").Value = -59819.332
